# Apply weekly Fruta/Hortaliza update for Pepino dulce - Macroferia Regional de Talca
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 23, pushing the existing row 23 down to row 24.
#    (The inserted row will carry row 23's old formatting, e.g. the date
#    style in column D.)
$ws.Rows("23").Insert()

# 2) Update row 21 with this week's new reading (date + volume).
$ws.Range("D21").Value = 44627
$ws.Range("J21").Value = 300

# 3) Update row 22 with the latest "Primera" quality price data, which
#    previously lived only in the now-shifted historical rows.
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 15000
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = 15000
$ws.Range("P22").Value = 833

# 4) Populate the newly inserted row 23 with the historical values that
#    used to sit in row 22 before this week's update.
$ws.Range("A23").Value = 5
$ws.Range("B23").Value = "Macroferia Regional de Talca"
$ws.Range("C23").Value = "Maule"
$ws.Range("D23").Value = 44396
$ws.Range("E23").Value = 7
$ws.Range("F23").Value = 100112043
$ws.Range("G23").Value = "Pepino dulce"
$ws.Range("H23").Value = "Cultivar IV Región"
$ws.Range("I23").Value = "Segunda"
$ws.Range("J23").Value = 150
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 12000
$ws.Range("M23").Value = 12000
$ws.Range("N23").Value = "$/bandeja 18 kilos"
$ws.Range("O23").Value = "Provincia de Limarí"
$ws.Range("P23").Value = 667
$ws.Range("Q23").Value = 18
$ws.Range("R23").Value = "Hortaliza"
